# 5th commit: update the test login data and move the cursor.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B2 held the shared string "test@9162" (with a mailto hyperlink left intact);
# extend it to "test@9162123".
$ws.Range("B2").Value = "test@9162123"

# Leave the selection where the author last left it when saving (B9 -> B17).
$ws.Range("B17").Select()
